$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1279
$ws.Range("E2").Value = 22
$ws.Range("F2").Value = 22
$ws.Range("G2").Value = -20
$ws.Range("H2").Value = -40
$ws.Range("I2").Value = -40
$ws.Range("K2").Value = 1448
$ws.Range("L2").Value = 899
$ws.Range("M2").Value = 549
$ws.Range("N2").Value = 549
$ws.Range("P2").Value = 78
$ws.Range("Q2").Value = -108
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 135
$ws.Range("T2").Value = 32
$ws.Range("U2").Value = -140
$ws.Range("V2").Value = 514
$ws.Range("W2").Value = 1.71
$ws.Range("X2").Value = -3.16
$ws.Range("Y2").Value = -7.49
$ws.Range("Z2").Value = -2.82
$ws.Range("AA2").Value = 163.71
$ws.Range("AB2").Value = 607.94
$ws.Range("AC2").Value = -1116
$ws.Range("AD2").Value = -19.08
$ws.Range("AE2").Value = 14185
$ws.Range("AF2").Value = 1.5
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 3871070
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 1160
$ws.Range("E3").Value = -274
$ws.Range("F3").Value = -245
$ws.Range("G3").Value = -329
$ws.Range("H3").Value = -367
$ws.Range("I3").Value = -367
$ws.Range("K3").Value = 1140
$ws.Range("L3").Value = 944
$ws.Range("M3").Value = 196
$ws.Range("N3").Value = 196
$ws.Range("P3").Value = 78
$ws.Range("Q3").Value = -79
$ws.Range("R3").Value = 20
$ws.Range("S3").Value = 55
$ws.Range("T3").Value = 24
$ws.Range("U3").Value = -103
$ws.Range("V3").Value = 537
$ws.Range("W3").Value = -23.64
$ws.Range("X3").Value = -31.65
$ws.Range("Y3").Value = -98.51
$ws.Range("Z3").Value = -28.36
$ws.Range("AA3").Value = 481.33
$ws.Range("AB3").Value = 151.21
$ws.Range("AC3").Value = -9456
$ws.Range("AD3").Value = -1.89
$ws.Range("AE3").Value = 5046
$ws.Range("AF3").Value = 3.54
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 3887446
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 667
$ws.Range("E4").Value = -326
$ws.Range("F4").Value = -325
$ws.Range("G4").Value = -396
$ws.Range("H4").Value = -392
$ws.Range("I4").Value = -392
$ws.Range("K4").Value = 767
$ws.Range("L4").Value = 533
$ws.Range("M4").Value = 234
$ws.Range("N4").Value = 234
$ws.Range("P4").Value = 307
$ws.Range("Q4").Value = -100
$ws.Range("R4").Value = 21
$ws.Range("S4").Value = 174
$ws.Range("T4").Value = 19
$ws.Range("U4").Value = -119
$ws.Range("V4").Value = 270
$ws.Range("W4").Value = -48.89
$ws.Range("X4").Value = -58.8
$ws.Range("Y4").Value = -182.44
$ws.Range("Z4").Value = -41.13
$ws.Range("AA4").Value = 228.03
$ws.Range("AB4").Value = -22.35
$ws.Range("AC4").Value = -8495
$ws.Range("AD4").Value = -0.49
$ws.Range("AE4").Value = 1902
$ws.Range("AF4").Value = 2.19
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 12288666
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 268
$ws.Range("E5").Value = -159
$ws.Range("F5").Value = -159
$ws.Range("G5").Value = -347
$ws.Range("H5").Value = -379
$ws.Range("I5").Value = -379
$ws.Range("K5").Value = 494
$ws.Range("L5").Value = 191
$ws.Range("M5").Value = 303
$ws.Range("N5").Value = 303
$ws.Range("P5").Value = 151
$ws.Range("Q5").Value = -85
$ws.Range("R5").Value = -9
$ws.Range("S5").Value = 41
$ws.Range("T5").Value = 11
$ws.Range("U5").Value = -95
$ws.Range("V5").Value = 62
$ws.Range("W5").Value = -59.4
$ws.Range("X5").Value = -141.39
$ws.Range("Y5").Value = -141.25
$ws.Range("Z5").Value = -60.16
$ws.Range("AA5").Value = 62.95
$ws.Range("AB5").Value = 103.92
$ws.Range("AC5").Value = -2680
$ws.Range("AD5").Value = -1.29
$ws.Range("AE5").Value = 1004
$ws.Range("AF5").Value = 3.45
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 30221846
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 474
$ws.Range("E6").Value = -88
$ws.Range("F6").Value = -88
$ws.Range("G6").Value = -138
$ws.Range("H6").Value = -140
$ws.Range("I6").Value = -135
$ws.Range("K6").Value = 1071
$ws.Range("L6").Value = 339
$ws.Range("M6").Value = 732
$ws.Range("N6").Value = 675
$ws.Range("P6").Value = 223
$ws.Range("Q6").Value = -69
$ws.Range("R6").Value = -543
$ws.Range("S6").Value = 682
$ws.Range("T6").Value = 11
$ws.Range("U6").Value = -80
$ws.Range("V6").Value = 208
$ws.Range("W6").Value = -18.48
$ws.Range("X6").Value = -29.59
$ws.Range("Y6").Value = -27.64
$ws.Range("Z6").Value = -17.91
$ws.Range("AA6").Value = 46.23
$ws.Range("AB6").Value = 202.53
$ws.Range("AC6").Value = -329
$ws.Range("AD6").Value = -5.79
$ws.Range("AE6").Value = 1515
$ws.Range("AF6").Value = 1.26
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 44549792

# Rows 7,8,9: clear all data columns except A,B,C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
